$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append Gm11633 to the COX (Oxidative phosphorylation complex IV) gene list
$ws.Range("B32").Value = "COX1, COX2, COX3, Cox4i1, Cox4i2, Cox5a, Cox5b, Cox6a1, Cox6a2, Cox6b1, Cox6b2, Cox6c, Cox7a1, Cox7a2, Cox7a2l, Cox7b, Cox7b2, Cox7c, Cox8a, Cox8b, Cox8c, Gm11633"

# Add Gm19340 prefix and Ndufab1-ps to the Nduf (complex I) gene list
$ws.Range("B33").Value = "Gm19340, Ndufa1, Ndufa10, Ndufa11, Ndufa12, Ndufa13, Ndufa2, Ndufa3, Ndufa4, Ndufa4l2, Ndufa5, Ndufa6, Ndufa7, Ndufa8, Ndufa9, Ndufab1, Ndufab1-ps, Ndufb1, Ndufb10, Ndufb11, Ndufb2, Ndufb3, Ndufb4, Ndufb4b, Ndufb4c, Ndufb5, Ndufb6, Ndufb7, Ndufb8, Ndufb9, Ndufc1, Ndufc2, Ndufs1, Ndufs2, Ndufs3, Ndufs4, Ndufs5, Ndufs6, Ndufs7, Ndufs8, Ndufv1, Ndufv2, Ndufv3"

# Insert Gm6293 into the Uqcr/CYTB (complex III) gene list
$ws.Range("B35").Value = "CYTB, Cyc1, Gm6293, Uqcr10, Uqcr11, Uqcrb, Uqcrc1, Uqcrc2, Uqcrfs1, Uqcrh, Uqcrq"
